$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Append a new sentence to the "RV_VILLAGE_NAME" bullet about the
#    SETTLEMENT field ordering.
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("The RV_VILLAGE_NAME field needs to be moved up in the tree, as its now asked after RV_HOUSEHOLD_NUMBER.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$appendRange = $d.Range($rng.End, $rng.End)
$appendRange.InsertAfter(" Also, the SETTLEMENT field is not in the right order on the form.")

# ---------------------------------------------------------------------
# 2) Remove the whole "Is there ever a case ... PROC RV_WARD is not
#    necessary." bullet -- it is no longer relevant.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("Is there ever a case where this program will be called without all of the IDs filled in? If not, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraToDelete = $rng2.Paragraphs(1)
$paraToDelete.Range.Delete()

# ---------------------------------------------------------------------
# 3) Rework the "demographic checks" bullet: reword the opening clause,
#    switch "and" separators to commas in the PROC list, tighten the
#    closing clause, and append a new sentence about the head of
#    household not needing to be on the first line.
# ---------------------------------------------------------------------

# "RV_B5 and PROC RV_B7" -> "RV_B5, PROC RV_B7"
$r1 = $d.Content
$r1.Find.Execute("RV_B5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterB5 = $d.Range($r1.End, $d.Content.End)
$afterB5.Find.Execute(" and ", $true, $false, $false, $false, $false, $true, 1, $false, ", ", 1) | Out-Null

# "RV_B7 and PROC DEMOGRAPHICS_ROSTER" -> "RV_B7, and PROC DEMOGRAPHICS_ROSTER"
$r2 = $d.Content
$r2.Find.Execute("RV_B7", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterB7 = $d.Range($r2.End, $d.Content.End)
$afterB7.Find.Execute(" and ", $true, $false, $false, $false, $false, $true, 1, $false, ", and ", 1) | Out-Null

# "Do you really want to have the demographic checks in" -> "The demographic checks in"
$r3 = $d.Content
$r3.Find.Execute("Do you really want to have the demographic checks in", $true, $false, $false, $false, $false, $true, 1, $false, "The demographic checks in", 2) | Out-Null

# "? If so, then they probably should match the checks in the household data entry application?"
#   -> " need to match the checks in the household data entry application? Remember that in this program the head does not have to be on the first line."
$r4 = $d.Content
$r4.Find.Execute("? If so, then they probably should match the checks in the household data entry application?", $true, $false, $false, $false, $false, $true, 1, $false, " need to match the checks in the household data entry application? Remember that in this program the head does not have to be on the first line.", 2) | Out-Null
